# The workbook schema table on Sheet1 renames a handful of column-name
# cells (column B) to match the updated database schema naming:
#   device_id        -> unique_id          (row 4,  devices table)
#   device_type      -> os_type            (row 5,  devices table)
#   preference_id    -> pref_id            (row 11, preferences table)
#   preference_name  -> pref_name          (row 12, preferences table)
#   preference_id    -> pref_id            (row 17, device_preference table)
#   preference_value -> pref_value         (row 18, device_preference table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B5").Value = "os_type"
$ws.Range("B4").Value = "unique_id"
$ws.Range("B11").Value = "pref_id"
$ws.Range("B12").Value = "pref_name"
$ws.Range("B17").Value = "pref_id"
$ws.Range("B18").Value = "pref_value"
